$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header values in row 1 for columns P and Q (matching header style of O1)
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15
$ws.Range("O1").Copy()
$ws.Range("P1:Q1").PasteSpecial(-4122)

# Update existing columns I, K, M, O for rows 2-25, and add new columns P, Q
for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 9).Value = 2   # column I
    $ws.Cells.Item($r, 11).Value = 1  # column K
    $ws.Cells.Item($r, 13).Value = 2  # column M
    $ws.Cells.Item($r, 15).Value = 1  # column O
    $ws.Cells.Item($r, 16).Value = 2  # column P
    $ws.Cells.Item($r, 17).Value = 2  # column Q
}
